{"js": "// The first paragraph originally reads \"# \" + \"Examplish Lexicon\" (two runs).\n// It becomes a single run \"# test\", and eight new paragraphs are inserted\n// right after it (the new command-processor test lines + a blank line +\n// a relocated \"# Examplish Lexicon\" heading) before the rest of the\n// document (which starts with the pre-existing blank paragraph) continues\n// unchanged.\nconst body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\nconst originalFirst = body.paragraphs.items[0];\n\n// Build the replacement first paragraph (\"# test\") ahead of the original,\n// then drop the original two-run paragraph. Inserting a brand-new\n// paragraph (rather than editing the old runs in place) keeps the normal\n// empty <w:rPr/> on the surviving run, matching how Word represents a\n// plain, unformatted run.\nconst newFirst = originalFirst.insertParagraph(\"# test\", Word.InsertLocation.before);\n\nconst newLines = [\n  \"\\\\pos vtest _{tense}\",\n  \"ka = vtest go\",\n  \"\\\\inflect vtest n = tense PRES\",\n  \"so = vtest eat\",\n  \"\\\\inflect vtest k = tense PAST\",\n  \"mi = vtest speak\",\n  \"\",\n  \"# Examplish Lexicon\",\n];\n\nlet anchor = newFirst;\nfor (const line of newLines) {\n  anchor = anchor.insertParagraph(line, Word.InsertLocation.after);\n}\n\noriginalFirst.delete();\n\nawait context.sync();\n", "ps1": "# The first paragraph originally reads \"# \" + \"Examplish Lexicon\" (two runs).\n# It becomes a single run \"# test\", and eight new paragraphs are inserted\n# right after it (the new command-processor test lines + a blank line +\n# a relocated \"# Examplish Lexicon\" heading) before the rest of the\n# document (which starts with the pre-existing blank paragraph) continues\n# unchanged.\n$d = $word.ActiveDocument\n\n$p1 = $d.Paragraphs(1)\n\n# Insert the replacement first paragraph (\"# test\") ahead of the original\n# two-run paragraph, then delete the original. Building a brand-new\n# paragraph (rather than editing the old runs in place) keeps the normal\n# empty <w:rPr/> on the surviving run, matching how Word represents a\n# plain, unformatted run.\n$p1.Range.InsertParagraphBefore()\n$d.Paragraphs(1).Range.Text = \"# test\"\n\n$newLines = @(\n  \"\\pos vtest _{tense}\",\n  \"ka = vtest go\",\n  \"\\inflect vtest n = tense PRES\",\n  \"so = vtest eat\",\n  \"\\inflect vtest k = tense PAST\",\n  \"mi = vtest speak\",\n  \"\",\n  \"# Examplish Lexicon\"\n)\n\n$idx = 1\nforeach ($line in $newLines) {\n  $anchor = $d.Paragraphs($idx)\n  $anchor.Range.InsertParagraphAfter()\n  $idx = $idx + 1\n  if ($line -ne \"\") {\n    $d.Paragraphs($idx).Range.Text = $line\n  }\n}\n\n# The original paragraph (\"# \" + \"Examplish Lexicon\") has now been pushed\n# down right after the newly-inserted block; remove it.\n$oldFirstIdx = $idx + 1\n$d.Paragraphs($oldFirstIdx).Range.Delete()\n"}
